$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(10,11,12,14,24,25,28,29,51,53,61,66,67,68,71,72,73,74,75,76,77)
foreach ($r in $rows) {
    $ws.Range("F$r").ClearContents()
}
